$d = $word.ActiveDocument

# --- 1. Insert "Upload first project:" paragraph before the existing "git init" line (paragraph 43) ---
$d.Paragraphs(43).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs(43).Range.Text = "Upload first project:"

# After this insertion the original git command paragraphs shifted down by one:
#   44 git init
#   45 git add .
#   46 git commit -m "first commit"
#   47 git branch -M main
#   48 git remote add origin <link>
#   49 git push -u origin main

# --- 2. Append a clean blank paragraph right after "git push -u origin main" (paragraph 49) ---
# Using InsertXML (rather than InsertParagraphAfter) keeps the new paragraph a true empty
# <w:p/> instead of inheriting an empty run/format from its source paragraph.
$xmlFrag49 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>git push -u origin main</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(49).Range.InsertXML($xmlFrag49) | Out-Null
# Paragraph 50 is now the new blank paragraph.

# --- 3. Add "Faced problem..." paragraph and the four follow-up git command paragraphs ---
$d.Paragraphs(50).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(51).Range.Text = "Faced problem, only root directory gets uploaded, run this:"

$d.Paragraphs(51).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(52).Range.Text = "git rm -r --cached ."

$d.Paragraphs(52).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(53).Range.Text = "git add ."

$d.Paragraphs(53).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(54).Range.Text = 'git commit -m "Re-add all files to Git"'

$d.Paragraphs(54).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(55).Range.Text = "git push origin main --force"

# --- 4. Indent all the git-command paragraphs (spacing after=0, left indent 720 twips / 36pt) ---
$targets = 44,45,46,47,48,49,52,53,54,55
foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    $p.Format.SpaceAfter = 0
    $p.Format.LeftIndent = 36
}

# --- 5. Rebuild paragraph 53 ("git add .") with the grammar-checked run split ---
# git <proofErr gramStart>add<proofErr gramEnd> .
$xmlFrag53 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">git </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>add</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(53).Range.InsertXML($xmlFrag53) | Out-Null
